$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (gdp_mmm_usd): new values for columns J2:AS2
$row2Values = @(
    659.0425967, 649.22546224, 639.40832778, 629.59119332, 619.77405886, 609.9569244,
    616.62913864, 623.30135288, 629.97356712, 636.64578136, 643.3179956, 666.67096334,
    690.02393108, 713.37689882, 736.72986656, 760.0828342999999, 784.39470634, 808.70657838,
    833.01845042, 857.33032246, 881.6421945, 906.0350418, 930.4278891, 954.8207364,
    979.2135837, 1003.606431, 1027.7438144, 1051.8811978, 1076.0185812, 1100.1559646,
    1124.293348, 1149.068216, 1173.843084, 1198.617952, 1223.39282, 1248.167688
)

$startCol = 10  # column J
for ($i = 0; $i -lt $row2Values.Length; $i++) {
    $ws.Cells.Item(2, $startCol + $i).Value = $row2Values[$i]
}

# Row 8 (elasticity_gnrl_rate_occupancy_to_gdppc): constant -0.1 across J8:AS8
for ($c = 10; $c -le 45; $c++) {
    $ws.Cells.Item(8, $c).Value = -0.1
}

# Row 9 (frac_gnrl_eating_red_meat): constant 1 across J9:AS9
for ($c = 10; $c -le 45; $c++) {
    $ws.Cells.Item(9, $c).Value = 1
}

# Row 13 (occrateinit_gnrl_occupancy): constant 3.145207224 across J13:AS13
for ($c = 10; $c -le 45; $c++) {
    $ws.Cells.Item(13, $c).Value = 3.145207224
}
